$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume table (Tue Apr  2 14:43:04 UTC 2024 refresh).
# Price (D) cells are forced to Text format before assignment so that
# numeric-looking strings (e.g. "1.00", "0.999") are preserved verbatim
# instead of being auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.075.71"
$ws.Range("E2").Value = "  -4.43%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.310.01"
$ws.Range("E3").Value = "  -5.94%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.49"
$ws.Range("E5").Value = "  -4.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.56"
$ws.Range("E6").Value = "  -5.65%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  -2.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.302.96"
$ws.Range("E9").Value = "  -5.77%  "

$ws.Range("E10").Value = "  -6.77%  "

$ws.Range("E11").Value = "  -4.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.75"
$ws.Range("E12").Value = "  -8.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000267"
$ws.Range("E13").Value = "  -6.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "635.45"
$ws.Range("E14").Value = "  -1.20%  "

$ws.Range("E15").Value = "  -6.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.840.02"
$ws.Range("E16").Value = "  -5.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.083.73"
$ws.Range("E17").Value = "  -4.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "17.98"
$ws.Range("E18").Value = "  -1.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.117"
$ws.Range("E19").Value = "  -3.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.309.23"
$ws.Range("E20").Value = "  -6.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.48"
$ws.Range("E21").Value = "  -7.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.905"
$ws.Range("E22").Value = "  -4.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.62"
$ws.Range("E23").Value = "  -1.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "106.86"
$ws.Range("E24").Value = "  +5.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.07"
$ws.Range("E25").Value = "  -6.74%  "

$ws.Range("E26").Value = "  -7.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.02"
$ws.Range("E27").Value = "  -0.38%  "

$ws.Range("E28").Value = "  -6.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.53"
$ws.Range("E29").Value = "  -5.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.77"
$ws.Range("E30").Value = "  -7.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.71"
$ws.Range("E31").Value = "  -6.29%  "

$ws.Range("E32").Value = "  -2.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.36"
$ws.Range("E33").Value = "  -5.23%  "

$ws.Range("E34").Value = "  -4.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "557.25"
$ws.Range("E35").Value = "  +10.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.105"
$ws.Range("E36").Value = "  -3.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.728.97"
$ws.Range("E37").Value = "  +0.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "57.44"
$ws.Range("E38").Value = "  -6.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.52"
$ws.Range("E40").Value = "  -1.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.74"
$ws.Range("E41").Value = "  -6.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0718"
$ws.Range("E42").Value = "  -9.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.128"
$ws.Range("E43").Value = "  -3.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("B44").Value = "CoreDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D44").Value = "3.36"
$ws.Range("E44").Value = "  +24.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "32.27"
$ws.Range("E45").Value = "  -6.53%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "0.342"
$ws.Range("E46").Value = "  -6.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0416"
$ws.Range("E47").Value = "  -5.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.24"
$ws.Range("E48").Value = "  -4.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.63"
$ws.Range("E49").Value = "  -7.19%  "

$ws.Range("E50").Value = "  -3.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -0.25%  "
